$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4097.946
$ws.Range("I15").Value = 4097.946
$ws.Range("K15").Value = 12293.838
$ws.Range("M15").Value = -12124.838

$ws.Range("H129").Value = 846.1975
$ws.Range("J129").Value = 883.72974
$ws.Range("L129").Value = 2651.18922
$ws.Range("N129").Value = -12651.18922

$ws.Range("H137").Value = 1090.579
$ws.Range("I137").Value = 1026.25
$ws.Range("J137").Value = 1270.7
$ws.Range("K137").Value = 3078.75
$ws.Range("L137").Value = 3812.1
$ws.Range("M137").Value = -528.75
$ws.Range("N137").Value = -8912.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 10750.1
$ws.Range("I2").Value = 832
$ws.Range("J2").Value = 100013
$ws.Range("K2").Value = 832
$ws.Range("L2").Value = 100013
$ws.Range("M2").Value = -719
$ws.Range("N2").Value = -100239

$ws.Range("H116").Value = 10750.1
$ws.Range("I116").Value = 832
$ws.Range("J116").Value = 100013
$ws.Range("K116").Value = 832
$ws.Range("L116").Value = 100013
$ws.Range("M116").Value = 1462
$ws.Range("N116").Value = -104601

$ws.Range("H122").Value = 1806.3572
$ws.Range("I122").Value = 1853.6364
$ws.Range("J122").Value = 1633
$ws.Range("K122").Value = 5560.9092
$ws.Range("L122").Value = 4899
$ws.Range("M122").Value = -3110.9092
$ws.Range("N122").Value = -9799

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 10750.1
$ws.Range("I3").Value = 832
$ws.Range("J3").Value = 100013
$ws.Range("K3").Value = 832
$ws.Range("L3").Value = 100013
$ws.Range("M3").Value = -718
$ws.Range("N3").Value = -100241

$ws.Range("H86").Value = 2717
$ws.Range("I86").Value = 2952.0667
$ws.Range("J86").Value = 2129.3333
$ws.Range("K86").Value = 2952.0667
$ws.Range("L86").Value = 2129.3333
$ws.Range("M86").Value = -1829.0667
$ws.Range("N86").Value = -4375.3333

$ws.Range("H89").Value = 2717
$ws.Range("I89").Value = 2952.0667
$ws.Range("J89").Value = 2129.3333
$ws.Range("K89").Value = 14760.3335
$ws.Range("L89").Value = 10646.6665
$ws.Range("M89").Value = -9144.333499999999
$ws.Range("N89").Value = -21878.6665

$ws.Range("H107").Value = 1628.8334
$ws.Range("I107").Value = 1203.3
$ws.Range("J107").Value = 3756.5
$ws.Range("K107").Value = 1203.3
$ws.Range("L107").Value = 3756.5
$ws.Range("M107").Value = 716.7
$ws.Range("N107").Value = -7596.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 10104
$ws.Range("I59").Value = 10104
$ws.Range("K59").Value = 10104
$ws.Range("M59").Value = -8959

$ws.Range("H86").Value = 3717637.5
$ws.Range("I86").Value = 6668986
$ws.Range("K86").Value = 6668986
$ws.Range("M86").Value = -6667863

$ws.Range("H89").Value = 3717637.5
$ws.Range("I89").Value = 6668986
$ws.Range("K89").Value = 33344930
$ws.Range("M89").Value = -33339314

$ws.Range("H107").Value = 557
$ws.Range("I107").Value = 444.2353
$ws.Range("J107").Value = 748.7
$ws.Range("K107").Value = 444.2353
$ws.Range("L107").Value = 748.7
$ws.Range("M107").Value = 1475.7647
$ws.Range("N107").Value = -4588.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 578.55554
$ws.Range("I5").Value = 578.55554
$ws.Range("K5").Value = 1735.66662
$ws.Range("M5").Value = -1623.66662

$ws.Range("H135").Value = 578.55554
$ws.Range("I135").Value = 578.55554
$ws.Range("K135").Value = 5206.99986
$ws.Range("M135").Value = -2671.99986

$ws.Range("H139").Value = 2519.0667
$ws.Range("I139").Value = 3209.5557
$ws.Range("J139").Value = 1483.3334
$ws.Range("K139").Value = 9628.667099999999
$ws.Range("L139").Value = 4450.0002
$ws.Range("M139").Value = -4488.667099999999
$ws.Range("N139").Value = -14730.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1666.6818
$ws.Range("I102").Value = 1817.0555
$ws.Range("J102").Value = 990
$ws.Range("K102").Value = 1817.0555
$ws.Range("L102").Value = 990
$ws.Range("M102").Value = -195.0554999999999
$ws.Range("N102").Value = -4234

$ws.Range("H107").Value = 1150.25
$ws.Range("I107").Value = 1200.3334
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1200.3334
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 719.6666
$ws.Range("N107").Value = -4840

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = ""

$ws.Range("H132").Value = 2943.25
$ws.Range("I132").Value = 2642
$ws.Range("K132").Value = 7926
$ws.Range("M132").Value = -5396

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2856.111
$ws.Range("I40").Value = 2800
$ws.Range("J40").Value = 2968.3333
$ws.Range("K40").Value = 2800
$ws.Range("L40").Value = 2968.3333
$ws.Range("M40").Value = -2664
$ws.Range("N40").Value = -3240.3333

$ws.Range("H122").Value = 25003228
$ws.Range("I122").Value = 35717612
$ws.Range("J122").Value = 2997
$ws.Range("K122").Value = 107152836
$ws.Range("L122").Value = 8991
$ws.Range("M122").Value = -107150386
$ws.Range("N122").Value = -13891

$ws.Range("H132").Value = 37611.965
$ws.Range("I132").Value = 1506.35
$ws.Range("J132").Value = 127876
$ws.Range("K132").Value = 4519.049999999999
$ws.Range("L132").Value = 383628
$ws.Range("M132").Value = -1989.049999999999
$ws.Range("N132").Value = -388688

$ws.Range("H136").Value = 13076.444
$ws.Range("I136").Value = 26546
$ws.Range("J136").Value = 2300.8
$ws.Range("K136").Value = 79638
$ws.Range("L136").Value = 6902.400000000001
$ws.Range("M136").Value = -77088
$ws.Range("N136").Value = -12002.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 535
$ws.Range("I113").Value = 358
$ws.Range("J113").Value = 1154.5
$ws.Range("K113").Value = 1074
$ws.Range("L113").Value = 3463.5
$ws.Range("M113").Value = 1096
$ws.Range("N113").Value = -7803.5

$ws.Range("H132").Value = 4831.3125
$ws.Range("I132").Value = 4392.8335
$ws.Range("J132").Value = 6146.75
$ws.Range("K132").Value = 13178.5005
$ws.Range("L132").Value = 18440.25
$ws.Range("M132").Value = -10648.5005
$ws.Range("N132").Value = -23500.25

$ws.Range("H136").Value = 1124.6666
$ws.Range("I136").Value = 949.7
$ws.Range("J136").Value = 1999.5
$ws.Range("K136").Value = 2849.1
$ws.Range("L136").Value = 5998.5
$ws.Range("M136").Value = -299.1000000000004
$ws.Range("N136").Value = -11098.5
